# "graphs for the paper"
# Re-run of the ContrastBasedAttack.xlsx "data" extraction: a header row
# (compress rate / attack name / epsilon / correct / counter /
#  correct rate (%) / time (sec)) that was missing before the second
# block of ContrastReductionAttack rows gets inserted at row 185 of the
# "data" sheet, pushing every row below it down by one. Formulas on the
# "results" sheet that reference data!F<row> shift automatically.
# View/selection state is also nudged to match where the author had
# scrolled to while building the new chart range.

$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1") selection moves down near the second chart ------
$wsSheet1 = $wb.Worksheets.Item("Sheet1")
$wsSheet1.Activate()
$wsSheet1.Range("C75").Select()

# --- data sheet: insert the missing header row at row 185 -------------
$wsData = $wb.Worksheets.Item("data")
$wsData.Activate()

$wsData.Rows.Item(185).Insert()

$wsData.Range("A185").Value = "compress rate"
$wsData.Range("B185").Value = " attack name"
$wsData.Range("C185").Value = " epsilon"
$wsData.Range("D185").Value = " correct"
$wsData.Range("E185").Value = " counter"
$wsData.Range("F185").Value = " correct rate (%)"
$wsData.Range("G185").Value = " time (sec)"

$wsData.Range("E185").Select()
